$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2; all existing rows (2..76) shift down to (3..77)
$ws.Rows.Item(2).Insert()

# The inserted row inherits formatting from the header row (row 1) by default.
# Reset it so it looks like a normal data row, then re-apply the date format
# used by every other "Fecha" (column D) cell.
$ws.Range("A2:T2").ClearFormats()
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate the new record (weekly refresh adds the most recent observation).
$ws.Range("A2").Value = 4
$ws.Range("B2").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C2").Value = "Los Lagos"
$ws.Range("D2").Value = 45282
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100103
$ws.Range("H2").Value = "Frutos de hueso (carozo)"
$ws.Range("I2").Value = 100103003
$ws.Range("J2").Value = "Damasco"
$ws.Range("K2").Value = "Castle Brite"
$ws.Range("L2").Value = "Tercera"
$ws.Range("M2").Value = 400
$ws.Range("N2").Value = 21000
$ws.Range("O2").Value = 21000
$ws.Range("P2").Value = 21000
$ws.Range("Q2").Value = "`$/caja 15 kilos"
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 1400
$ws.Range("T2").Value = 15
